# Update "想去人数" (F column) values on the 展览 and 全部类型 sheets
# to reflect the latest scrape snapshot, per commit 456a3b4.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F6").Value = 548
$wsExhibit.Range("F7").Value = 1696
$wsExhibit.Range("F8").Value = 27
$wsExhibit.Range("F11").Value = 1699
$wsExhibit.Range("F13").Value = 90
$wsExhibit.Range("F14").Value = 413
$wsExhibit.Range("F18").Value = 29
$wsExhibit.Range("F21").Value = 483
$wsExhibit.Range("F24").Value = 232
$wsExhibit.Range("F25").Value = 251

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 548
$wsAll.Range("F7").Value = 1696
$wsAll.Range("F9").Value = 27
$wsAll.Range("F12").Value = 1699
$wsAll.Range("F14").Value = 90
$wsAll.Range("F15").Value = 413
$wsAll.Range("F19").Value = 29
$wsAll.Range("F22").Value = 484
$wsAll.Range("F25").Value = 232
$wsAll.Range("F26").Value = 251

$wb.Save()
